$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.994.45"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.167.68"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.02%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "627.59"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +33.42%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.372"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.166.86"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.15%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +14.13%  "
$ws.Range("E11").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.60%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.76"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.91%  "
$ws.Range("E13").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "35.13"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +7.74%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.801.83"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.757.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.47%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.189.55"
$ws.Range("D18").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +12.65%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.75"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +9.41%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "479.03"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +12.25%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.15%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +10.36%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.24"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.56%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.90"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +16.63%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.93"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +10.26%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.41"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.50%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.343.55"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.88%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.40%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +23.50%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -20.33%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +42.70%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "525.94"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.27%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.58%  "
$ws.Range("E36").ClearFormats()
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Kaspa"
$ws.Range("B37").ClearFormats()
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C37").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.145"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +8.95%  "
$ws.Range("E37").ClearFormats()
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("B38").ClearFormats()
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C38").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.63"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("E38").ClearFormats()
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "RenderToken"
$ws.Range("B39").ClearFormats()
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("C39").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.98"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.24%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.61%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0905"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +29.91%  "
$ws.Range("E41").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.423"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +17.61%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("E44").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.53%  "
$ws.Range("E45").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +19.43%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "150.34"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.00%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.62"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +9.18%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +12.56%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.38"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.08%  "
$ws.Range("E51").ClearFormats()
